$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, using the same bold/bordered style as the
# other header cells (match formatting of the existing G1 header cell).
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.Item(7).LineStyle = 1
$ws.Range("H1").Borders.Item(8).LineStyle = 1
$ws.Range("H1").Borders.Item(9).LineStyle = 1
$ws.Range("H1").Borders.Item(10).LineStyle = 1

# Map each subject name (column A) to its diagnosis label: Control -> 0, MDD -> 1
$labelMap = @{
    "Control 26" = 0
    "Control 33" = 0
    "Control 36" = 0
    "Control 49" = 0
    "Control 2"  = 0
    "MDD 37"     = 1
    "MDD 24"     = 1
    "MDD 6"      = 1
    "MDD 54"     = 1
    "MDD 21"     = 1
}

for ($r = 2; $r -le 21; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 8).Value = $labelMap[$name]
}

# A handful of prediction/error/cross-entropy values were refreshed from a
# re-fit of the model (higher precision refit results).
$ws.Range("D4").Value  = 0.3694740740422722
$ws.Range("E4").Value  = 0.3694740740422722

$ws.Range("D5").Value  = 0.2972392597445135
$ws.Range("E5").Value  = 0.2972392597445135

$ws.Range("D9").Value  = 0.3640526047991846
$ws.Range("E9").Value  = 0.6359473952008154

$ws.Range("D11").Value = 0.3894742636485105
$ws.Range("E11").Value = 0.6105257363514895
$ws.Range("F11").Value = 0.6206362247467041
